$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.712.91"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "3.047.35"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.374"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "3.574.40"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000162"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "57.744.31"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.48%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.041.98"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.498"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "0.0₃0895"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "3.093.12"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.652"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "2.275.13"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0255"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.10%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.934"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.732"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "252.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.79%  "
